$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "308.16"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.43%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.05"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.05%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.041"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.41%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07648"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-2.71%"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-1.87%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.616"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-4.23%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.459"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-3.91%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9082"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.73%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1012"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-5.34%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1774"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.44%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09220"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2.02%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-4.25%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1052"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.82%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001258"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.58%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005808"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.70%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.61%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.3270"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.43%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.776"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-5.99%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1356"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.76%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04160"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.01%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001219"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-2.10%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004086"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.42%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "5.96%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003007"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.12%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02411"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-1.06%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05179"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-2.29%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007787"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.74%"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-3.55%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007082"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-6.22%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001968"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.54%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007476"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-8.74%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3055"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.69%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006376"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-6.22%"
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.006177"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "80.02%"
$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004400"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "6.43%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002000"
